$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the existing row 3 ("The Blue Toes"), shifting
# it (and the rows below it) down by one.
$ws.Rows.Item(3).Insert()

# Fill in the new team's data in the freshly-inserted row 3.
$ws.Range("A3").Value = "Footloose"
$ws.Range("B3").Value = "Eindhoven"
$ws.Range("C3").Value = "NTDS_Eindhoven.xlsx"

# Match the final selection recorded in the saved file.
$ws.Range("C4").Select()
